# StructureDefinition-communication-vendor.xlsx
# Update FHIR StructureDefinition metadata (Metadata sheet) and refresh the
# Short/Definition columns for the root Extension element (Elements sheet)
# to match the new 6.0.0 publication of the IG.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Metadata" sheet (Property / Value)
$ws2 = $wb.Worksheets.Item(2)   # "Elements" sheet

# --- Metadata sheet -------------------------------------------------------

# Version bump: 5.0.0 -> 6.0.0
$ws1.Range("B3").Value = "6.0.0"

# Date updated
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$ws1.Range("B9").Value = "Alvearie Team"

# The old "Contact" row is replaced with a "Jurisdiction" row
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Remove the now-duplicate second "Contact" / "No display for ContactDetail" row
$ws1.Range("A11").EntireRow.Delete()

# --- Elements sheet --------------------------------------------------------

# Row 2 is the root Extension element; its Short/Definition columns picked up
# the resource-specific text instead of the generic placeholder.
$ws2.Range("K2").Value = "Communication Vendor"
$ws2.Range("L2").Value = "Vendor used to send the communication"
